# B6-PowerPoint.pptx edit
#
# 1) Three tables (slides 14, 15, 16) get their table style switched from
#    the deck's custom "Table_0" style ({A6EC6344-3417-4519-9ADC-33E58843A832})
#    to the built-in theme style {58A87983-26CA-4992-A5B1-665382C0C60D}
#    (PowerPoint's "Medium Style 2 - Accent 1").
#
# 2) The presentation's theme colour scheme is swapped from the "Integral /
#    Red Violet" palette to the default "Office" palette (i.e. the deck's
#    Design swaps from the pink/violet Integral theme to the plain Office
#    Theme colours).

$p = $ppt.ActivePresentation

# --- 1) Re-style the three summary tables -------------------------------
$newTableStyle = "{58A87983-26CA-4992-A5B1-665382C0C60D}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Swap the theme colour scheme (Integral -> Office) ---------------
function RGBValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office theme colours, in the same dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink
# order exposed by ThemeColorScheme.Colors(1..12).
$officeColors = @(
    (RGBValue 0x00 0x00 0x00),  # dk1
    (RGBValue 0xFF 0xFF 0xFF),  # lt1
    (RGBValue 0x44 0x54 0x6A),  # dk2
    (RGBValue 0xE7 0xE6 0xE6),  # lt2
    (RGBValue 0x5B 0x9B 0xD5),  # accent1
    (RGBValue 0xED 0x7D 0x31),  # accent2
    (RGBValue 0xA5 0xA5 0xA5),  # accent3
    (RGBValue 0xFF 0xC0 0x00),  # accent4
    (RGBValue 0x44 0x72 0xC4),  # accent5
    (RGBValue 0x70 0xAD 0x47),  # accent6
    (RGBValue 0x05 0x63 0xC1),  # hlink
    (RGBValue 0x95 0x4F 0x72)   # folHlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
